$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.599.34'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '2.605.09'
$ws.Range("E3").Value = '  +1.68%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +3.68%  '
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("E10").Value = '  +2.07%  '
$ws.Range("E11").Value = '  +1.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.130'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.96%  '
$ws.Range("D13").Value = '3.054.31'
$ws.Range("E13").Value = '  +1.05%  '
$ws.Range("D14").Value = '60.632.55'
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000141'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.47%  '
$ws.Range("D17").Value = '2.614.09'
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '356.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.89%  '
$ws.Range("E20").Value = '  +3.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.46%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("E23").Value = '  +1.74%  '
$ws.Range("E24").Value = '  +2.92%  '
$ws.Range("E25").Value = '  +0.69%  '
$ws.Range("D26").Value = '2.718.18'
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.996'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.18%  '
$ws.Range("D28").Value = '0.0₃0844'
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.37'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.50'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.91%  '
$ws.Range("E33").Value = '  +2.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.67'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.949'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.83%  '
$ws.Range("E37").Value = '  +1.27%  '
$ws.Range("E38").Value = '  +2.40%  '
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.35'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.845'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '288.53'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.50%  '
$ws.Range("E43").Value = '  +2.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.627'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.40%  '
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.996'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0237'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.71%  '
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("D51").Value = '1.982.84'
$ws.Range("E51").Value = '  -1.40%  '
